$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 421
$ws.Cells.Item(421, 1).Value = 420
$ws.Cells.Item(421, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(421, 3).Value = "2:20 PM"
$ws.Cells.Item(421, 4).Value = "FR9280"
$ws.Cells.Item(421, 5).Value = "Agadir"
$ws.Cells.Item(421, 6).Value = "(AGA)"
$ws.Cells.Item(421, 7).Value = "Ryanair "
$ws.Cells.Item(421, 8).Value = "B738"
$ws.Cells.Item(421, 9).Value = "(EI-EVA)"
$ws.Cells.Item(421, 10).Value = "2:33 PM"
$ws.Cells.Item(421, 12).Value = "0 hours, 13 minutes"

# Row 422
$ws.Cells.Item(422, 1).Value = 421
$ws.Cells.Item(422, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(422, 3).Value = "2:30 PM"
$ws.Cells.Item(422, 4).Value = "OS598"
$ws.Cells.Item(422, 5).Value = "Vienna"
$ws.Cells.Item(422, 6).Value = "(VIE)"
$ws.Cells.Item(422, 7).Value = "Austrian Airlines "
$ws.Cells.Item(422, 8).Value = "E195"
$ws.Cells.Item(422, 9).Value = "(OE-LWA)"
$ws.Cells.Item(422, 10).Value = "2:30 PM"
$ws.Cells.Item(422, 12).Value = "0 hours, 0 minutes"

# Row 423
$ws.Cells.Item(423, 1).Value = 422
$ws.Cells.Item(423, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(423, 3).Value = "2:50 PM"
$ws.Cells.Item(423, 4).Value = "FR6355"
$ws.Cells.Item(423, 5).Value = "Alicante"
$ws.Cells.Item(423, 6).Value = "(ALC)"
$ws.Cells.Item(423, 7).Value = "Buzz "
$ws.Cells.Item(423, 8).Value = "B38M"
$ws.Cells.Item(423, 9).Value = "(SP-RZB)"
$ws.Cells.Item(423, 10).Value = "2:59 PM"
$ws.Cells.Item(423, 12).Value = "0 hours, 9 minutes"

# Row 424
$ws.Cells.Item(424, 1).Value = 423
$ws.Cells.Item(424, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(424, 3).Value = "3:00 PM"
$ws.Cells.Item(424, 4).Value = "LO3906"
$ws.Cells.Item(424, 5).Value = "Warsaw"
$ws.Cells.Item(424, 6).Value = "(WAW)"
$ws.Cells.Item(424, 7).Value = "LOT "
$ws.Cells.Item(424, 8).Value = "E195"
$ws.Cells.Item(424, 9).Value = "(SP-LNP)"
$ws.Cells.Item(424, 10).Value = "3:09 PM"
$ws.Cells.Item(424, 12).Value = "0 hours, 9 minutes"

# Row 425
$ws.Cells.Item(425, 1).Value = 424
$ws.Cells.Item(425, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(425, 3).Value = "4:20 PM"
$ws.Cells.Item(425, 4).Value = "FR6314"
$ws.Cells.Item(425, 5).Value = "Marseille"
$ws.Cells.Item(425, 6).Value = "(MRS)"
$ws.Cells.Item(425, 7).Value = "Ryanair "
$ws.Cells.Item(425, 8).Value = "B38M"
$ws.Cells.Item(425, 9).Value = "(9H-VUW)"
$ws.Cells.Item(425, 10).Value = "4:34 PM"
$ws.Cells.Item(425, 12).Value = "0 hours, 14 minutes"

# Row 426
$ws.Cells.Item(426, 1).Value = 425
$ws.Cells.Item(426, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(426, 3).Value = "4:25 PM"
$ws.Cells.Item(426, 4).Value = "FR6214"
$ws.Cells.Item(426, 5).Value = "Dortmund"
$ws.Cells.Item(426, 6).Value = "(DTM)"
$ws.Cells.Item(426, 7).Value = "Buzz "
$ws.Cells.Item(426, 8).Value = "B38M"
$ws.Cells.Item(426, 9).Value = "(SP-RZF)"
$ws.Cells.Item(426, 10).Value = "4:29 PM"
$ws.Cells.Item(426, 12).Value = "0 hours, 4 minutes"

# Row 427
$ws.Cells.Item(427, 1).Value = 426
$ws.Cells.Item(427, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(427, 3).Value = "4:45 PM"
$ws.Cells.Item(427, 4).Value = "LO3507"
$ws.Cells.Item(427, 5).Value = "Olsztyn"
$ws.Cells.Item(427, 6).Value = "(SZY)"
$ws.Cells.Item(427, 7).Value = "LOT (Star Alliance Livery) "
$ws.Cells.Item(427, 8).Value = "E75S"
$ws.Cells.Item(427, 9).Value = "(SP-LIO)"
$ws.Cells.Item(427, 10).Value = "4:39 PM"
$ws.Cells.Item(427, 12).Value = "0 hours, -6 minutes"

# Row 428
$ws.Cells.Item(428, 1).Value = 427
$ws.Cells.Item(428, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(428, 3).Value = "4:50 PM"
$ws.Cells.Item(428, 4).Value = "KL1996"
$ws.Cells.Item(428, 5).Value = "Amsterdam"
$ws.Cells.Item(428, 6).Value = "(AMS)"
$ws.Cells.Item(428, 7).Value = "KLM "
$ws.Cells.Item(428, 8).Value = "E190"
$ws.Cells.Item(428, 9).Value = "(PH-EZR)"
$ws.Cells.Item(428, 10).Value = "5:00 PM"
$ws.Cells.Item(428, 12).Value = "0 hours, 10 minutes"

# Row 429
$ws.Cells.Item(429, 1).Value = 428
$ws.Cells.Item(429, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(429, 3).Value = "4:55 PM"
$ws.Cells.Item(429, 4).Value = "W65051"
$ws.Cells.Item(429, 5).Value = "Larnaca"
$ws.Cells.Item(429, 6).Value = "(LCA)"
$ws.Cells.Item(429, 7).Value = "Wizz Air "
$ws.Cells.Item(429, 8).Value = "A21N"
$ws.Cells.Item(429, 9).Value = "(HA-LZI)"
$ws.Cells.Item(429, 10).Value = "5:06 PM"
$ws.Cells.Item(429, 12).Value = "0 hours, 11 minutes"

# Row 430
$ws.Cells.Item(430, 1).Value = 429
$ws.Cells.Item(430, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(430, 3).Value = "5:00 PM"
$ws.Cells.Item(430, 4).Value = "FR239"
$ws.Cells.Item(430, 5).Value = "Gdansk"
$ws.Cells.Item(430, 6).Value = "(GDN)"
$ws.Cells.Item(430, 7).Value = "Ryanair "
$ws.Cells.Item(430, 8).Value = "B738"
$ws.Cells.Item(430, 9).Value = "(SP-RSO)"
$ws.Cells.Item(430, 10).Value = "5:15 PM"
$ws.Cells.Item(430, 12).Value = "0 hours, 15 minutes"

# Row 431
$ws.Cells.Item(431, 1).Value = 430
$ws.Cells.Item(431, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(431, 3).Value = "5:15 PM"
$ws.Cells.Item(431, 4).Value = "FR4627"
$ws.Cells.Item(431, 5).Value = "Stockholm"
$ws.Cells.Item(431, 6).Value = "(ARN)"
$ws.Cells.Item(431, 7).Value = "Ryanair "
$ws.Cells.Item(431, 8).Value = "B38M"
$ws.Cells.Item(431, 9).Value = "(9H-VUJ)"
$ws.Cells.Item(431, 10).Value = "5:46 PM"
$ws.Cells.Item(431, 12).Value = "0 hours, 31 minutes"

# Row 432
$ws.Cells.Item(432, 1).Value = 431
$ws.Cells.Item(432, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(432, 3).Value = "5:30 PM"
$ws.Cells.Item(432, 4).Value = "FR5892"
$ws.Cells.Item(432, 5).Value = "Eindhoven"
$ws.Cells.Item(432, 6).Value = "(EIN)"
$ws.Cells.Item(432, 7).Value = "Ryanair "
$ws.Cells.Item(432, 8).Value = "B738"
$ws.Cells.Item(432, 9).Value = "(SP-RSH)"
$ws.Cells.Item(432, 10).Value = "5:48 PM"
$ws.Cells.Item(432, 12).Value = "0 hours, 18 minutes"

# Row 433
$ws.Cells.Item(433, 1).Value = 432
$ws.Cells.Item(433, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(433, 3).Value = "5:35 PM"
$ws.Cells.Item(433, 4).Value = "FR6222"
$ws.Cells.Item(433, 5).Value = "Bournemouth"
$ws.Cells.Item(433, 6).Value = "(BOH)"
$ws.Cells.Item(433, 7).Value = "Ryanair "
$ws.Cells.Item(433, 8).Value = "B738"
$ws.Cells.Item(433, 9).Value = "(EI-EGB)"
$ws.Cells.Item(433, 10).Value = "5:50 PM"
$ws.Cells.Item(433, 12).Value = "0 hours, 15 minutes"

# Row 434
$ws.Cells.Item(434, 1).Value = 433
$ws.Cells.Item(434, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(434, 3).Value = "5:40 PM"
$ws.Cells.Item(434, 4).Value = "DY1029"
$ws.Cells.Item(434, 5).Value = "Bergen"
$ws.Cells.Item(434, 6).Value = "(BGO)"
$ws.Cells.Item(434, 7).Value = "Norwegian "
$ws.Cells.Item(434, 8).Value = "B738"
$ws.Cells.Item(434, 9).Value = "(LN-NIH)"
$ws.Cells.Item(434, 10).Value = "5:40 PM"
$ws.Cells.Item(434, 12).Value = "0 hours, 0 minutes"

# Row 435
$ws.Cells.Item(435, 1).Value = 434
$ws.Cells.Item(435, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(435, 3).Value = "5:45 PM"
$ws.Cells.Item(435, 4).Value = "W65003"
$ws.Cells.Item(435, 5).Value = "London"
$ws.Cells.Item(435, 6).Value = "(LTN)"
$ws.Cells.Item(435, 7).Value = "Wizz Air "
$ws.Cells.Item(435, 8).Value = "A21N"
$ws.Cells.Item(435, 9).Value = "(HA-LVH)"
$ws.Cells.Item(435, 10).Value = "6:02 PM"
$ws.Cells.Item(435, 12).Value = "0 hours, 17 minutes"

# Row 436
$ws.Cells.Item(436, 1).Value = 435
$ws.Cells.Item(436, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(436, 3).Value = "5:55 PM"
$ws.Cells.Item(436, 4).Value = "W65041"
$ws.Cells.Item(436, 5).Value = "Bergen"
$ws.Cells.Item(436, 6).Value = "(BGO)"
$ws.Cells.Item(436, 7).Value = "Wizz Air "
$ws.Cells.Item(436, 8).Value = "A321"
$ws.Cells.Item(436, 9).Value = "(HA-LXO)"
$ws.Cells.Item(436, 10).Value = "6:39 PM"
$ws.Cells.Item(436, 12).Value = "0 hours, 44 minutes"

# Row 437
$ws.Cells.Item(437, 1).Value = 436
$ws.Cells.Item(437, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(437, 3).Value = "6:00 PM"
$ws.Cells.Item(437, 4).Value = "FR8672"
$ws.Cells.Item(437, 5).Value = "Girona"
$ws.Cells.Item(437, 6).Value = "(GRO)"
$ws.Cells.Item(437, 7).Value = "Ryanair "
$ws.Cells.Item(437, 8).Value = "B738"
$ws.Cells.Item(437, 9).Value = "(SP-RSA)"
$ws.Cells.Item(437, 10).Value = "7:14 PM"
$ws.Cells.Item(437, 12).Value = "1 hours, 14 minutes"

# Row 438
$ws.Cells.Item(438, 1).Value = 437
$ws.Cells.Item(438, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(438, 3).Value = "6:05 PM"
$ws.Cells.Item(438, 4).Value = "FR2363"
$ws.Cells.Item(438, 5).Value = "London"
$ws.Cells.Item(438, 6).Value = "(STN)"
$ws.Cells.Item(438, 7).Value = "Ryanair "
$ws.Cells.Item(438, 8).Value = "B38M"
$ws.Cells.Item(438, 9).Value = "(SP-RZL)"
$ws.Cells.Item(438, 10).Value = "6:27 PM"
$ws.Cells.Item(438, 12).Value = "0 hours, 22 minutes"

# Row 439
$ws.Cells.Item(439, 1).Value = 438
$ws.Cells.Item(439, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(439, 3).Value = "6:15 PM"
$ws.Cells.Item(439, 4).Value = "LO3924"
$ws.Cells.Item(439, 5).Value = "Warsaw"
$ws.Cells.Item(439, 6).Value = "(WAW)"
$ws.Cells.Item(439, 7).Value = "LOT "
$ws.Cells.Item(439, 8).Value = "E190"
$ws.Cells.Item(439, 9).Value = "(SP-LMF)"
$ws.Cells.Item(439, 10).Value = "6:22 PM"
$ws.Cells.Item(439, 12).Value = "0 hours, 7 minutes"

# Row 440
$ws.Cells.Item(440, 1).Value = 439
$ws.Cells.Item(440, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(440, 3).Value = "6:35 PM"
$ws.Cells.Item(440, 4).Value = "FR3722"
$ws.Cells.Item(440, 5).Value = "Billund"
$ws.Cells.Item(440, 6).Value = "(BLL)"
$ws.Cells.Item(440, 7).Value = "Ryanair "
$ws.Cells.Item(440, 8).Value = "B738"
$ws.Cells.Item(440, 9).Value = "(9H-QCY)"
$ws.Cells.Item(440, 10).Value = "6:41 PM"
$ws.Cells.Item(440, 12).Value = "0 hours, 6 minutes"

# Row 441
$ws.Cells.Item(441, 1).Value = 440
$ws.Cells.Item(441, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(441, 3).Value = "6:35 PM"
$ws.Cells.Item(441, 4).Value = "U26526"
$ws.Cells.Item(441, 5).Value = "London"
$ws.Cells.Item(441, 6).Value = "(LGW)"
$ws.Cells.Item(441, 7).Value = "easyJet "
$ws.Cells.Item(441, 8).Value = "A320"
$ws.Cells.Item(441, 9).Value = "(G-EZTR)"
$ws.Cells.Item(441, 10).Value = "6:56 PM"
$ws.Cells.Item(441, 12).Value = "0 hours, 21 minutes"

# Row 442
$ws.Cells.Item(442, 1).Value = 441
$ws.Cells.Item(442, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(442, 3).Value = "7:25 PM"
$ws.Cells.Item(442, 4).Value = "W65077"
$ws.Cells.Item(442, 5).Value = "Stockholm"
$ws.Cells.Item(442, 6).Value = "(NYO)"
$ws.Cells.Item(442, 7).Value = "Wizz Air "
$ws.Cells.Item(442, 8).Value = "A21N"
$ws.Cells.Item(442, 9).Value = "(HA-LVO)"
$ws.Cells.Item(442, 10).Value = "7:27 PM"
$ws.Cells.Item(442, 12).Value = "0 hours, 2 minutes"

# Row 443
$ws.Cells.Item(443, 1).Value = 442
$ws.Cells.Item(443, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(443, 3).Value = "7:40 PM"
$ws.Cells.Item(443, 4).Value = "FR9663"
$ws.Cells.Item(443, 5).Value = "Rome"
$ws.Cells.Item(443, 6).Value = "(CIA)"
$ws.Cells.Item(443, 7).Value = "Ryanair "
$ws.Cells.Item(443, 8).Value = "B738"
$ws.Cells.Item(443, 9).Value = "(9H-QCO)"
$ws.Cells.Item(443, 10).Value = "7:48 PM"
$ws.Cells.Item(443, 12).Value = "0 hours, 8 minutes"

# Row 444
$ws.Cells.Item(444, 1).Value = 443
$ws.Cells.Item(444, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(444, 3).Value = "7:45 PM"
$ws.Cells.Item(444, 4).Value = "BA873"
$ws.Cells.Item(444, 5).Value = "London"
$ws.Cells.Item(444, 6).Value = "(LHR)"
$ws.Cells.Item(444, 7).Value = "British Airways "
$ws.Cells.Item(444, 8).Value = "A320"
$ws.Cells.Item(444, 9).Value = "(G-EUYA)"
$ws.Cells.Item(444, 10).Value = "8:11 PM"
$ws.Cells.Item(444, 12).Value = "0 hours, 26 minutes"

# Row 445
$ws.Cells.Item(445, 1).Value = 444
$ws.Cells.Item(445, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(445, 3).Value = "7:45 PM"
$ws.Cells.Item(445, 4).Value = "FR7115"
$ws.Cells.Item(445, 5).Value = "Malta"
$ws.Cells.Item(445, 6).Value = "(MLA)"
$ws.Cells.Item(445, 7).Value = "Ryanair "
$ws.Cells.Item(445, 8).Value = "B738"
$ws.Cells.Item(445, 9).Value = "(9H-QBG)"
$ws.Cells.Item(445, 10).Value = "8:00 PM"
$ws.Cells.Item(445, 12).Value = "0 hours, 15 minutes"

# Row 446
$ws.Cells.Item(446, 1).Value = 445
$ws.Cells.Item(446, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(446, 3).Value = "7:55 PM"
$ws.Cells.Item(446, 4).Value = "U23818"
$ws.Cells.Item(446, 5).Value = "Paris"
$ws.Cells.Item(446, 6).Value = "(CDG)"
$ws.Cells.Item(446, 7).Value = "easyJet "
$ws.Cells.Item(446, 8).Value = "A320"
$ws.Cells.Item(446, 9).Value = "(OE-IZF)"
$ws.Cells.Item(446, 10).Value = "8:13 PM"
$ws.Cells.Item(446, 12).Value = "0 hours, 18 minutes"

# Row 447
$ws.Cells.Item(447, 1).Value = 446
$ws.Cells.Item(447, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(447, 3).Value = "8:15 PM"
$ws.Cells.Item(447, 4).Value = "FR5724"
$ws.Cells.Item(447, 5).Value = "Nottingham"
$ws.Cells.Item(447, 6).Value = "(EMA)"
$ws.Cells.Item(447, 7).Value = "Ryanair "
$ws.Cells.Item(447, 8).Value = "B738"
$ws.Cells.Item(447, 9).Value = "(EI-EVK)"
$ws.Cells.Item(447, 10).Value = "8:40 PM"
$ws.Cells.Item(447, 12).Value = "0 hours, 25 minutes"

# Row 448
$ws.Cells.Item(448, 1).Value = 447
$ws.Cells.Item(448, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(448, 3).Value = "8:30 PM"
$ws.Cells.Item(448, 4).Value = "FR6361"
$ws.Cells.Item(448, 5).Value = "Shannon"
$ws.Cells.Item(448, 6).Value = "(SNN)"
$ws.Cells.Item(448, 7).Value = "Ryanair "
$ws.Cells.Item(448, 8).Value = "B738"
$ws.Cells.Item(448, 9).Value = "(EI-EKH)"
$ws.Cells.Item(448, 10).Value = "8:38 PM"
$ws.Cells.Item(448, 12).Value = "0 hours, 8 minutes"

# Row 449
$ws.Cells.Item(449, 1).Value = 448
$ws.Cells.Item(449, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(449, 3).Value = "8:35 PM"
$ws.Cells.Item(449, 4).Value = "FR3365"
$ws.Cells.Item(449, 5).Value = "Berlin"
$ws.Cells.Item(449, 6).Value = "(BER)"
$ws.Cells.Item(449, 7).Value = "Buzz "
$ws.Cells.Item(449, 8).Value = "B38M"
$ws.Cells.Item(449, 9).Value = "(SP-RZF)"
$ws.Cells.Item(449, 10).Value = "8:50 PM"
$ws.Cells.Item(449, 12).Value = "0 hours, 15 minutes"

# Row 450
$ws.Cells.Item(450, 1).Value = 449
$ws.Cells.Item(450, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(450, 3).Value = "9:25 PM"
$ws.Cells.Item(450, 4).Value = "FR5396"
$ws.Cells.Item(450, 5).Value = "Dublin"
$ws.Cells.Item(450, 6).Value = "(DUB)"
$ws.Cells.Item(450, 7).Value = "Ryanair "
$ws.Cells.Item(450, 8).Value = "B738"
$ws.Cells.Item(450, 9).Value = "(EI-DCJ)"
$ws.Cells.Item(450, 10).Value = "9:27 PM"
$ws.Cells.Item(450, 12).Value = "0 hours, 2 minutes"

# Row 451
$ws.Cells.Item(451, 1).Value = 450
$ws.Cells.Item(451, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(451, 3).Value = "9:25 PM"
$ws.Cells.Item(451, 4).Value = "LO3922"
$ws.Cells.Item(451, 5).Value = "Warsaw"
$ws.Cells.Item(451, 6).Value = "(WAW)"
$ws.Cells.Item(451, 7).Value = "LOT (Star Alliance Livery) "
$ws.Cells.Item(451, 8).Value = "E170"
$ws.Cells.Item(451, 9).Value = "(SP-LDK)"
$ws.Cells.Item(451, 10).Value = "9:29 PM"
$ws.Cells.Item(451, 12).Value = "0 hours, 4 minutes"

# Row 452
$ws.Cells.Item(452, 1).Value = 451
$ws.Cells.Item(452, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(452, 3).Value = "9:40 PM"
$ws.Cells.Item(452, 4).Value = "DY1043"
$ws.Cells.Item(452, 5).Value = "Oslo"
$ws.Cells.Item(452, 6).Value = "(OSL)"
$ws.Cells.Item(452, 7).Value = "Norwegian "
$ws.Cells.Item(452, 8).Value = "B738"
$ws.Cells.Item(452, 9).Value = "(LN-NII)"
$ws.Cells.Item(452, 10).Value = "9:35 PM"
$ws.Cells.Item(452, 12).Value = "0 hours, -5 minutes"
